$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Table S1 - Plasticity AIC")
$ws.Range("G2").Value = 1.2926268490263
$ws.Range("H2").Value = 0.649718186195778
$ws.Range("E3").Value = 281.9
$ws.Range("G3").Value = 1.2899740213025
$ws.Range("H3").Value = 0.653586411465873
$ws.Range("G4").Value = 1.33685374011254
$ws.Range("H4").Value = 0.660209766997358
$ws.Range("G5").Value = 1.32703425546414
$ws.Range("H5").Value = 0.660312676146787
$ws.Range("G6").Value = 1.25915911436328
$ws.Range("H6").Value = 0.657644659640462
$ws.Range("F7").Value = 0.1382
$ws.Range("G7").Value = 1.4151360221094
$ws.Range("H7").Value = 0.683221885809563
$ws.Range("G8").Value = 1.20315338344266
$ws.Range("H8").Value = 0.610559911143328
$ws.Range("G9").Value = 1.30997747698609
$ws.Range("H9").Value = 0.617460634216792
$ws.Range("G10").Value = 1.37957299867158
$ws.Range("H10").Value = 0.631205876260462
$ws.Range("G11").Value = 1.29607785162445
$ws.Range("H11").Value = 0.629451678765452
$ws.Range("G12").Value = 1.37596637996336
$ws.Range("H12").Value = 0.634439131239188
$ws.Range("G13").Value = 1.34632437229469
$ws.Range("H13").Value = 0.634659628138949
$ws.Range("G14").Value = 1.16376490169629
$ws.Range("H14").Value = 0.564285238534006
$ws.Range("G15").Value = 1.17631329356072
$ws.Range("H15").Value = 0.547207956382125
$ws.Range("G16").Value = 1.1637156190192
$ws.Range("H16").Value = 0.570253769164272
$ws.Range("E17").Value = 202.69
$ws.Range("G17").Value = 1.17633511372074
$ws.Range("H17").Value = 0.552401500246739
$ws.Range("D18").Value = 192.02
$ws.Range("G18").Value = 1.17064356267687
$ws.Range("H18").Value = 0.561821146309639
$ws.Range("G19").Value = 1.15904635034672
$ws.Range("H19").Value = 0.587675913316996

$ws = $wb.Worksheets.Item("Table S2 - PERMANOVA")
$ws.Range("C2").Value = 61072
$ws.Range("D2").Value = 0.208
$ws.Range("E2").Value = 8.15
$ws.Range("C3").Value = 7471
$ws.Range("D3").Value = 0.025
$ws.Range("E3").Value = 2.99
$ws.Range("F3").Value = 0.08728
$ws.Range("C4").Value = 24705
$ws.Range("D4").Value = 0.084
$ws.Range("E4").Value = 9.89
$ws.Range("F4").Value = 0.002
$ws.Range("C5").Value = 199740
$ws.Range("D5").Value = 0.682
$ws.Range("C6").Value = 292988
$ws.Range("C7").Value = 97850
$ws.Range("D7").Value = 0.087
$ws.Range("E7").Value = 14.29
$ws.Range("C8").Value = 26676
$ws.Range("D8").Value = 0.024
$ws.Range("E8").Value = 1.3
$ws.Range("F8").Value = 0.31246
$ws.Range("C9").Value = 519372
$ws.Range("D9").Value = 0.46
$ws.Range("E9").Value = 75.84
$ws.Range("C10").Value = 486202
$ws.Range("D10").Value = 0.43
$ws.Range("C11").Value = 1130099
$ws.Range("C12").Value = 157
$ws.Range("E12").Value = 0.11
$ws.Range("F12").Value = 0.76016
$ws.Range("C13").Value = 25414
$ws.Range("D13").Value = 0.18
$ws.Range("E13").Value = 18.47
$ws.Range("C14").Value = 30537
$ws.Range("D14").Value = 0.216
$ws.Range("E14").Value = 7.4
$ws.Range("F14").Value = 0.00133
$ws.Range("C15").Value = 85309
$ws.Range("D15").Value = 0.603
$ws.Range("C16").Value = 141417

$ws = $wb.Worksheets.Item("Table S3 - Plasticity GLM")
$ws.Range("D2").Value = 8.95
$ws.Range("E13").Value = 0.507
$ws.Range("E14").Value = 0.761

$ws = $wb.Worksheets.Item("Table S5 - Species PERMANOVA")
$ws.Columns.Item(5).ColumnWidth = 5.83   # OOXML width 5.71 -> 6.71
$ws.Range("C2").Value = 149393
$ws.Range("D2").Value = 0.038
$ws.Range("E2").Value = 8.24
$ws.Range("C3").Value = 17313
$ws.Range("D3").Value = 0.004
$ws.Range("E3").Value = 2.87
$ws.Range("F3").Value = 0.09194
$ws.Range("C4").Value = 58058
$ws.Range("D4").Value = 0.015
$ws.Range("E4").Value = 9.61
$ws.Range("F4").Value = 0.004
$ws.Range("C5").Value = 1642613
$ws.Range("D5").Value = 0.423
$ws.Range("E5").Value = 135.91
$ws.Range("C6").Value = 553351
$ws.Range("D6").Value = 0.143
$ws.Range("E6").Value = 45.78
$ws.Range("C7").Value = 90865
$ws.Range("D7").Value = 0.023
$ws.Range("E7").Value = 2.51
$ws.Range("F7").Value = 0.01732
$ws.Range("C8").Value = 77259
$ws.Range("D8").Value = 0.02
$ws.Range("E8").Value = 6.39
$ws.Range("F8").Value = 0.00466
$ws.Range("C9").Value = 1293204
$ws.Range("D9").Value = 0.333
$ws.Range("C10").Value = 3882055

$ws = $wb.Worksheets.Item("Table S6 - HostVsymb PERMANOVA")
$ws.Range("H2").Value = 61056
$ws.Range("I2").Value = 0.208
$ws.Range("J2").Value = 8.15
$ws.Range("H3").Value = 7468
$ws.Range("I3").Value = 0.025
$ws.Range("J3").Value = 2.99
$ws.Range("K3").Value = 0.08461
$ws.Range("H4").Value = 24705
$ws.Range("I4").Value = 0.084
$ws.Range("J4").Value = 9.9
$ws.Range("K4").Value = 0.002
$ws.Range("H5").Value = 199684
$ws.Range("I5").Value = 0.682
$ws.Range("H6").Value = 292913
$ws.Range("H7").Value = 26899
$ws.Range("I7").Value = 0.024
$ws.Range("J7").Value = 1.31
$ws.Range("K7").Value = 0.27781
$ws.Range("H8").Value = 515173
$ws.Range("I8").Value = 0.456
$ws.Range("J8").Value = 75.24
$ws.Range("H9").Value = 101793
$ws.Range("I9").Value = 0.09
$ws.Range("J9").Value = 14.87
$ws.Range("K9").Value = 0.00067
$ws.Range("H10").Value = 486140
$ws.Range("I10").Value = 0.43
$ws.Range("H11").Value = 1130005
$ws.Range("H12").Value = 29037
$ws.Range("I12").Value = 0.205
$ws.Range("J12").Value = 7.04
$ws.Range("H13").Value = 26338
$ws.Range("I13").Value = 0.186
$ws.Range("J13").Value = 19.15
$ws.Range("H14").Value = 724
$ws.Range("I14").Value = 0.005
$ws.Range("J14").Value = 0.53
$ws.Range("K14").Value = 0.48168
$ws.Range("H15").Value = 85288
$ws.Range("I15").Value = 0.603
$ws.Range("H16").Value = 141387
